$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at the top of the data (row 2),
# pushing every existing data row down by one and adding one new
# row at the very bottom of the range (former last row just shifts down).
$ws.Rows("2:2").Insert()

# Reset formatting on the freshly inserted row so it matches the plain
# (unstyled) look of the other data rows instead of inheriting the
# header row's bold/border/centered style.
$ws.Range("A2:R2").Style = "Normal"

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 'Macroferia Regional de Talca'
$ws.Range("C2").Value = 'Maule'
$ws.Range("D2").Value = 44473
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = 'Espárragos'
$ws.Range("H2").Value = 'Verde'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1200
$ws.Range("N2").Value = '$/kilo'
$ws.Range("O2").Value = 'Provincia de Linares'
$ws.Range("P2").Value = 1200
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'Hortaliza'
